# Updated cryptos list on Sat Mar 30 11:57:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B / C / E updates: plain text values (coin names, links, and
# percentage strings that already contain non-numeric characters such as
# "%" and surrounding spaces), so Excel will not reinterpret them as
# numbers and a direct .Value assignment keeps them as text.
$textUpdates = @{
    "E2"  = "  -0.16%  "
    "E3"  = "  +0.20%  "
    "E4"  = "  -0.10%  "
    "E5"  = "  -1.73%  "
    "E6"  = "  +6.45%  "
    "E7"  = "  -0.39%  "
    "E8"  = "  -0.08%  "
    "E9"  = "  -1.72%  "
    "E10" = "  +1.09%  "
    "E11" = "  +0.88%  "
    "E12" = "  -1.85%  "
    "E13" = "  +0.48%  "
    "E14" = "  -0.02%  "
    "E15" = "  -4.85%  "
    "B16" = "Chainlink"
    "C16" = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
    "E16" = "  +1.28%  "
    "B17" = "WrappedBTC"
    "C17" = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
    "E17" = "  +0.00%  "
    "E18" = "  -1.89%  "
    "E19" = "  -0.75%  "
    "E20" = "  +0.60%  "
    "E21" = "  +0.34%  "
    "E22" = "  +4.30%  "
    "E23" = "  +8.39%  "
    "E24" = "  +0.30%  "
    "E25" = "  -1.95%  "
    "E26" = "  +2.90%  "
    "E27" = "  -0.28%  "
    "E28" = "  +1.76%  "
    "E29" = "  -2.72%  "
    "E30" = "  +23.58%  "
    "E31" = "  +1.46%  "
    "E32" = "  +3.38%  "
    "E33" = "  +0.97%  "
    "E34" = "  -0.59%  "
    "E35" = "  +7.93%  "
    "E36" = "  +6.74%  "
    "E37" = "  -3.73%  "
    "E38" = "  -0.07%  "
    "E39" = "  +3.37%  "
    "E40" = "  -1.12%  "
    "E41" = "  -0.61%  "
    "E42" = "  -6.12%  "
    "E43" = "  -0.76%  "
    "E45" = "  -3.35%  "
    "E46" = "  -3.21%  "
    "E47" = "  -1.13%  "
    "E48" = "  +0.29%  "
    "E49" = "  -4.95%  "
    "E50" = "  +2.43%  "
    "E51" = "  +11.96%  "
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Column D updates: these are "Price" values that are formatted as plain
# text in the workbook (periods used as thousands separators, leading
# zeros, trailing zeros, etc.). Several of the new values look like
# ordinary numbers (e.g. "604.11", "37.00", "0.000247") and Excel would
# silently convert them to numeric cells - losing the exact textual
# representation (trailing zeros, fixed-point notation). Force the cell
# to Text format before writing the value, then clear the temporary
# formatting again so the cell ends up without an explicit style, just
# like the rest of the sheet.
$priceUpdates = @{
    "D2"  = "70.011.67"
    "D3"  = "3.553.03"
    "D5"  = "604.11"
    "D6"  = "197.54"
    "D9"  = "0.211"
    "D11" = "54.28"
    "D13" = "9.60"
    "D14" = "4.110.04"
    "D15" = "601.11"
    "D16" = "19.19"
    "D17" = "70.181.13"
    "D19" = "3.547.08"
    "D21" = "0.999"
    "D22" = "18.36"
    "D23" = "5.31"
    "D24" = "103.40"
    "D25" = "4.64"
    "D29" = "33.69"
    "D30" = "4.55"
    "D31" = "7.16"
    "D32" = "12.78"
    "D34" = "63.50"
    "D35" = "0.0₃0844"
    "D36" = "3.744.91"
    "D37" = "3.12"
    "D40" = "0.397"
    "D41" = "37.00"
    "D42" = "497.58"
    "D45" = "2.85"
    "D47" = "3.33"
    "D50" = "0.000247"
    "D51" = "1.31"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}
